$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.465.42"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.738.62"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.08"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4558"
$ws.Range("E7").Value = "  +8.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3529"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07400"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.28"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.39"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.040"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "1.739.28"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.28"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001053"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06340"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.55"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.718"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "27.506.64"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.081"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.87"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.02"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "1.937.24"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.61"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.044"
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09093"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.386"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02267"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.59"
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05947"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2055"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6226"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.374"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.693"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.03"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.705"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5784"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.97"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.927"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06840"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.09"
